$d = $word.ActiveDocument

# 1. "Search for address: " + "1DF F197 CFA0 " -> merge into a single run
$d.Content.Find.Execute("Search for address: 1DF F197 CFA0 ", $true, $false, $false, $false, $false, $true, 1, $false, "Search for address: 1DF F197 CFA0 ", 2)

# 2. Append note to "1DFF197B9A0 = Cash Manager - 10 X" (same run)
$d.Content.Find.Execute("1DFF197B9A0 = Cash Manager - 10 X", $true, $false, $false, $false, $false, $true, 1, $false, "1DFF197B9A0 = Cash Manager - 10 X -- There can be multiple addresses pointing to konfuze", 2)

# 3. "Search for address: " + "1DF F197 B990" -> merge into a single run
$d.Content.Find.Execute("Search for address: 1DF F197 B990", $true, $false, $false, $false, $false, $true, 1, $false, "Search for address: 1DF F197 B990", 2)

# 4. "1DFDC1DB030 = Cash Managers - 30 X" -> add trailing space, then append a NEW separate run with the note
$r4 = $d.Content
$r4.Find.Execute("1DFDC1DB030 = Cash Managers - 30 X", $true, $false, $false, $false, $false, $true, 1, $false, "1DFDC1DB030 = Cash Managers - 30 X ", 2)
$r4.Collapse(0)
$insertedStart = $r4.End
$r4.InsertAfter("-- There can be multiple addresses pointing to Cash Manager")
$insertedEnd = $r4.End
$newRun4 = $d.Range($insertedStart, $insertedEnd)
$newRun4.Font.Bold = $true
$newRun4.Font.Bold = $false

# 5. "Search for address: 1DF DC1D B000  " + "-- Following can have " -> merge into a single run
$d.Content.Find.Execute("Search for address: 1DF DC1D B000  -- Following can have ", $true, $false, $false, $false, $false, $true, 1, $false, "Search for address: 1DF DC1D B000  -- Following can have ", 2)

# 6. "Search for address: " + "1DFF1474300 " -> merge into a single run
$d.Content.Find.Execute("Search for address: 1DFF1474300 ", $true, $false, $false, $false, $false, $true, 1, $false, "Search for address: 1DFF1474300 ", 2)
